$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C18: 40 -> 40.33 (Distancia Media)
$ws.Range("C18").Value = 40.33

# Update C19: 4500 -> "4500.33" stored as text/shared string (Peso Total)
# Force text storage (so it round-trips as a shared string, like the target),
# then restore the original numeric-cell formatting that C19 had before.
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "4500.33"
$ws.Range("A19").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection / zoom on the sheet view
$ws.Range("C19").Select() | Out-Null
$excel.ActiveWindow.Zoom = 145
